$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "text"
$ws.Range("B1").Value = "extracted_model"
$ws.Range("C1").Value = "reason"

# Data row values
$ws.Range("A2").Value = "sansui led 24 jsv24nshd"
$ws.Range("C2").Value = "No matching spans found"

# B2 stays present-but-empty: touching a border (no-op at default) forces the
# cell to be materialized without pulling in a new style.
$ws.Range("B2").Borders.LineStyle = 0

# Build the header style (bold, centered/top aligned, thin box border) on a
# scratch cell first, then copy/paste-special just the formatting onto the
# header range in a single atomic style application, so only ONE new style
# combination is minted instead of one per property assignment.
$scratch = $ws.Range("E1")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160

$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()
